# Apply the lookup-table edit: replace the discrete -5.817e-4 coefficient
# grid (A1:H14) with a flat -0.05 placeholder, in preparation for an
# interpolation tool for smooth transitions between discrete values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:H14").Value = -0.05

# Reflect the updated selection (whole table selected) left behind by the edit.
$ws.Range("A1:H14").Select() | Out-Null
